# Updated the design files for the 9 port panel to use plastic parts with
# better isolation of stimulus LED light between ports.
#
# The "M3 hex spacer 3mm" mounting-hardware row (row 17, with its HARWIN /
# R30-9400300 supplier info) is no longer needed, so remove the whole row.
# The rows below it (M3 button head screw 12mm, M3 nut, M3 threaded insert,
# M3 Button head screw 6mm, the Electronic section header, and Speaker)
# shift up by one row automatically when the row is deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "M3 hex spacer 3mm" row entirely (shifts rows below up).
$ws.Rows(17).Delete()

# Restore the active-cell selection to where the author left it afterwards.
$ws.Range("B25").Select()
